$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "Thu Dec 07 21:34:43 EST 2023",
    "Thu Dec 07 21:34:52 EST 2023",
    "Thu Dec 07 21:35:01 EST 2023",
    "Thu Dec 07 21:35:10 EST 2023",
    "Thu Dec 07 21:35:19 EST 2023",
    "Thu Dec 07 21:35:29 EST 2023",
    "Thu Dec 07 21:35:38 EST 2023",
    "Thu Dec 07 21:35:47 EST 2023",
    "Thu Dec 07 21:35:56 EST 2023",
    "Thu Dec 07 21:36:05 EST 2023",
    "Thu Dec 07 21:36:14 EST 2023",
    "Thu Dec 07 21:36:23 EST 2023",
    "Thu Dec 07 21:36:32 EST 2023",
    "Thu Dec 07 21:36:41 EST 2023",
    "Thu Dec 07 21:36:50 EST 2023",
    "Thu Dec 07 21:36:59 EST 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
